$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Subregion" column header and blank placeholder values for
# the existing data rows (mirrors the Region normalization / sub_region
# support added in this revision).
$ws.Range("E1").Value = "Subregion"
$ws.Range("E2").Value = ""
$ws.Range("E3").Value = ""
